$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("I40").Value = 1849.75
$ws.Range("J40").Value = 2000
$ws.Range("K40").Value = 1849.75
$ws.Range("L40").Value = 2000
$ws.Range("M40").Value = -1674.75
$ws.Range("N40").Value = -2350
$ws.Range("H55").Value = 799.2222
$ws.Range("I55").Value = 531.6667
$ws.Range("J55").Value = 933
$ws.Range("K55").Value = 531.6667
$ws.Range("L55").Value = 933
$ws.Range("M55").Value = -317.6667
$ws.Range("N55").Value = -1361
$ws.Range("H135").Value = 1451.2
$ws.Range("I135").Value = 1144.8572
$ws.Range("J135").Value = 2166
$ws.Range("K135").Value = 10303.7148
$ws.Range("L135").Value = 19494
$ws.Range("M135").Value = -7768.7148
$ws.Range("N135").Value = -24564
$ws.Range("H138").Value = 6272.527
$ws.Range("J138").Value = 5913.527
$ws.Range("L138").Value = 17740.581
$ws.Range("N138").Value = -28020.581
$ws = $wb.Worksheets.Item(2)
$ws.Range("H2").Value = 2005.5834
$ws.Range("I2").Value = 2005.5834
$ws.Range("K2").Value = 2005.5834
$ws.Range("M2").Value = -1892.5834
$ws.Range("H5").Value = 200
$ws.Range("I5").Value = 200
$ws.Range("K5").Value = 200
$ws.Range("M5").Value = -88
$ws.Range("H32").Value = 24766.541
$ws.Range("J32").Value = 30186.227
$ws.Range("L32").Value = 30186.227
$ws.Range("N32").Value = -30760.227
$ws.Range("H61").Value = 2001.0834
$ws.Range("I61").Value = 1728.4546
$ws.Range("K61").Value = 1728.4546
$ws.Range("M61").Value = -1516.4546
$ws.Range("H74").Value = 6583.727
$ws.Range("I74").Value = 3331.6667
$ws.Range("J74").Value = 7803.25
$ws.Range("K74").Value = 3331.6667
$ws.Range("L74").Value = 7803.25
$ws.Range("M74").Value = -2457.6667
$ws.Range("N74").Value = -9551.25
$ws.Range("H77").Value = 6583.727
$ws.Range("I77").Value = 3331.6667
$ws.Range("J77").Value = 7803.25
$ws.Range("K77").Value = 16658.3335
$ws.Range("L77").Value = 39016.25
$ws.Range("M77").Value = -12290.3335
$ws.Range("N77").Value = -47752.25
$ws.Range("H88").Value = 5750
$ws.Range("I88").Value = 1500
$ws.Range("K88").Value = 1500
$ws.Range("M88").Value = -1094
$ws.Range("H91").Value = 5750
$ws.Range("I91").Value = 1500
$ws.Range("K91").Value = 1500
$ws.Range("M91").Value = -96
$ws.Range("H102").Value = 2168.4119
$ws.Range("I102").Value = 2388.5
$ws.Range("J102").Value = 2100.6924
$ws.Range("K102").Value = 2388.5
$ws.Range("L102").Value = 2100.6924
$ws.Range("M102").Value = -766.5
$ws.Range("N102").Value = -5344.6924
$ws.Range("H116").Value = 2005.5834
$ws.Range("I116").Value = 2005.5834
$ws.Range("K116").Value = 2005.5834
$ws.Range("M116").Value = 288.4166
$ws.Range("H122").Value = 913455.25
$ws.Range("H136").Value = 2001.0834
$ws.Range("I136").Value = 1728.4546
$ws.Range("K136").Value = 5185.3638
$ws.Range("M136").Value = -2635.3638
$ws = $wb.Worksheets.Item(3)
$ws.Range("H3").Value = 2005.5834
$ws.Range("I3").Value = 2005.5834
$ws.Range("K3").Value = 2005.5834
$ws.Range("M3").Value = -1891.5834
$ws.Range("H4").Value = 200
$ws.Range("I4").Value = 200
$ws.Range("K4").Value = 200
$ws.Range("M4").Value = -85
$ws.Range("H22").Value = 525.5
$ws.Range("I22").Value = 525.5
$ws.Range("K22").Value = 525.5
$ws.Range("M22").Value = -352.5
$ws.Range("H59").Value = 120500
$ws.Range("J59").Value = 120500
$ws.Range("L59").Value = 120500
$ws.Range("N59").Value = -122194
$ws.Range("H107").Value = 3992.04
$ws.Range("I107").Value = 2750.3684
$ws.Range("K107").Value = 2750.3684
$ws.Range("M107").Value = -830.3683999999998
$ws.Range("H134").Value = 4217.316
$ws.Range("I134").Value = 1914.125
$ws.Range("K134").Value = 5742.375
$ws.Range("M134").Value = -3207.375
$ws = $wb.Worksheets.Item(4)
$ws.Range("H31").Value = 4066.1396
$ws.Range("I31").Value = 2863.56
$ws.Range("K31").Value = 2863.56
$ws.Range("M31").Value = -2568.56
$ws.Range("H34").Value = 4066.1396
$ws.Range("I34").Value = 2863.56
$ws.Range("K34").Value = 2863.56
$ws.Range("M34").Value = -2661.56
$ws.Range("H58").Value = 8430.444
$ws.Range("I58").Value = 7387
$ws.Range("J58").Value = 8728.571
$ws.Range("K58").Value = 7387
$ws.Range("L58").Value = 8728.571
$ws.Range("M58").Value = -7184
$ws.Range("N58").Value = -9134.571
$ws.Range("H107").Value = 858.9
$ws.Range("I107").Value = 264.14285
$ws.Range("K107").Value = 264.14285
$ws.Range("M107").Value = 1655.85715
$ws.Range("H122").Value = 2000
$ws.Range("I122").Value = 2000
$ws.Range("K122").Value = 6000
$ws.Range("M122").Value = -3550
$ws.Range("H132").Value = 2056.1292
$ws.Range("I132").Value = 1975.08
$ws.Range("K132").Value = 5925.24
$ws.Range("M132").Value = -3395.24
$ws.Range("H134").Value = 3795.2222
$ws.Range("I134").Value = 3424.3333
$ws.Range("J134").Value = 4166.1113
$ws.Range("K134").Value = 10272.9999
$ws.Range("L134").Value = 12498.3339
$ws.Range("M134").Value = -7737.999899999999
$ws.Range("N134").Value = -17568.3339
$ws.Range("H136").Value = 8430.444
$ws.Range("I136").Value = 7387
$ws.Range("J136").Value = 8728.571
$ws.Range("K136").Value = 22161
$ws.Range("L136").Value = 26185.713
$ws.Range("M136").Value = -19611
$ws.Range("N136").Value = -31285.713
$ws = $wb.Worksheets.Item(5)
$ws.Range("H68").Value = 3254.8572
$ws.Range("I68").Value = 3698
$ws.Range("K68").Value = 11094
$ws.Range("M68").Value = -10283
$ws.Range("H71").Value = 3254.8572
$ws.Range("I71").Value = 3698
$ws.Range("K71").Value = 33282
$ws.Range("M71").Value = -29226
$ws.Range("H131").Value = 1404.0264
$ws.Range("J131").Value = 1570.2858
$ws.Range("L131").Value = 4710.857400000001
$ws.Range("N131").Value = -14790.8574
$ws.Range("H140").Value = 4682.3076
$ws.Range("I140").Value = 3715.5454
$ws.Range("K140").Value = 11146.6362
$ws.Range("M140").Value = -5966.636200000001
$ws = $wb.Worksheets.Item(6)
$ws.Range("H107").Value = 142.14285
$ws.Range("I107").Value = 150.83333
$ws.Range("J107").Value = 90
$ws.Range("K107").Value = 150.83333
$ws.Range("L107").Value = 90
$ws.Range("M107").Value = 1769.16667
$ws.Range("N107").Value = -3930
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
$ws.Range("H126").Value = 5000
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 15000
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -19940
$ws.Range("H132").Value = 2936.4517
$ws.Range("I132").Value = 2035.28
$ws.Range("K132").Value = 6105.84
$ws.Range("M132").Value = -3575.84
$ws = $wb.Worksheets.Item(7)
$ws.Range("H16").Value = 1547.5264
$ws.Range("I16").Value = 1523.7059
$ws.Range("K16").Value = 1523.7059
$ws.Range("M16").Value = -1353.7059
$ws.Range("H46").Value = 3625
$ws.Range("I46").Value = 2642.8572
$ws.Range("K46").Value = 2642.8572
$ws.Range("M46").Value = -2454.8572
$ws.Range("H93").Value = 1184.2858
$ws.Range("I93").Value = 622.5
$ws.Range("J93").Value = 1933.3334
$ws.Range("K93").Value = 622.5
$ws.Range("L93").Value = 1933.3334
$ws.Range("M93").Value = 625.5
$ws.Range("N93").Value = -4429.3334
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()
$ws = $wb.Worksheets.Item(8)
$ws.Range("H54").Value = 15071.429
$ws.Range("I54").Value = 1000
$ws.Range("K54").Value = 1000
$ws.Range("M54").Value = -480
$ws.Range("H81").Value = 3625.125
$ws.Range("I81").Value = 3250
$ws.Range("K81").Value = 6500
$ws.Range("M81").Value = -5439
$ws.Range("H84").Value = 3625.125
$ws.Range("I84").Value = 3250
$ws.Range("K84").Value = 32500
$ws.Range("M84").Value = -27196
$ws.Range("H113").Value = 1109.3684
$ws.Range("I113").Value = 552.3077
$ws.Range("J113").Value = 2316.3333
$ws.Range("K113").Value = 1656.9231
$ws.Range("L113").Value = 6948.999899999999
$ws.Range("M113").Value = 513.0769
$ws.Range("N113").Value = -11288.9999
$ws.Range("H132").Value = 2795.625
$ws.Range("I132").Value = 2248.6365
$ws.Range("J132").Value = 3999
$ws.Range("K132").Value = 6745.9095
$ws.Range("L132").Value = 11997
$ws.Range("M132").Value = -4215.9095
$ws.Range("N132").Value = -17057
$ws.Range("H136").Value = 204058.6
$ws.Range("I136").Value = 1152
$ws.Range("K136").Value = 3456
$ws.Range("M136").Value = -906
